$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# A2: "Version: ..." cell
$aboutSheet.Range("A2").Value = "Version: " + $newVersion

# A6: Recommended Citation cell
$citationNew = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for No. 4 Coal Mine (AL), United States, M1397, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"
$aboutSheet.Range("A6").Value = $citationNew

# S2:S27 on the data sheet hold the build_version string
for ($r = 2; $r -le 27; $r++) {
    $cell = $dataSheet.Cells.Item($r, 19)
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
